# Apply the edits described by the diff:
# 1. Update the date string "07152023" -> "07302023" wherever it appears in
#    the sheet (it is shared between the EffectiveDate and PreviousExpDate
#    columns on rows 2 and 3: F2, I2, F3, I3).
# 2. Update the sheet view: scroll the window so column J is the top-left
#    visible column, and change the active selection to G10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Replace the old date value with the new one, wherever it is used ---
$oldDate = "07152023"
$newDate = "07302023"

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value = $newDate
    }
}

# --- 2. Update the view: scroll position + active selection ---
# Scroll so column J becomes the left-most visible column …
$excel.ActiveWindow.ScrollColumn = 10
# … and move the active selection to G10.
$ws.Range("G10").Select()
